# Function to standardize data
# - Rename "Häufigkeit Blinken" header to "Häufigkeit Blinzeln (/min)"
# - Convert raw blink counts (column T) into a standardized rate by
#   multiplying each existing value by 15/29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header
$headerCell = $ws.Range("T1")
if ($headerCell.Value2 -eq "Häufigkeit Blinken") {
    $headerCell.Value2 = "Häufigkeit Blinzeln (/min)"
}

$factor = 15 / 29

# Standardize every populated value in column T (data rows start at row 2)
$lastRow = $ws.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 20)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val * $factor
    }
}
